$wb = $excel.ActiveWorkbook

# Delete row 6 on sheet R1 (the LTH0330 / SCECO+STB / Good row that was removed).
$wsR1 = $wb.Worksheets.Item("R1")
$wsR1.Rows.Item(6).Delete()

# Bump the "Elapsed Duration(Hrs)" (column G) values that were refreshed in the edit.
$wsR1.Range("G2").Value = "3929:44:11"
$wsR1.Range("G3").Value = "69:16:49"

$wsR2 = $wb.Worksheets.Item("R2")
$wsR2.Range("G2").Value = "12111:07:45"
$wsR2.Range("G3").Value = "3240:51:14"
$wsR2.Range("G4").Value = "479:02:48"

$wsR4 = $wb.Worksheets.Item("R4")
$wsR4.Range("G2").Value = "2956:57:34"
$wsR4.Range("G3").Value = "184:09:49"
$wsR4.Range("G4").Value = "72:22:14"
$wsR4.Range("G5").Value = "69:59:47"

$wsR5 = $wb.Worksheets.Item("R5")
$wsR5.Range("G2").Value = "430:56:33"

$wsR6 = $wb.Worksheets.Item("R6")
$wsR6.Range("G2").Value = "71:28:51"
